$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content (old No/Name/Data sample data)
$ws.Cells.Clear()

# Headers (row 1)
$headers = @("product_id", "category_id", "tax_id", "title", "description", "price", "language", "specifications", "is_discount", "is_highlight")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Font.Size = 10
}

# Data rows 2-3, filled column by column so shared-string insertion order
# matches the source workbook (title strings before language/spec strings).
$ws.Cells.Item(2, 1).Value = 3
$ws.Cells.Item(3, 1).Value = 4

$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(3, 2).Value = 1

$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(3, 3).Value = 1

$ws.Cells.Item(2, 4).Value = "Máy chấm công mới"
$ws.Cells.Item(3, 4).Value = "Máy chấm công mới A"

$ws.Cells.Item(2, 5).Value = "Máy chấm công mới"
$ws.Cells.Item(3, 5).Value = "Máy chấm công mới"

$ws.Cells.Item(2, 6).Value = 200000
$ws.Cells.Item(3, 6).Value = 400000

$ws.Cells.Item(2, 7).Value = "a"
$ws.Cells.Item(3, 7).Value = "a"

$ws.Cells.Item(2, 8).Value = "a"
$ws.Cells.Item(3, 8).Value = "a"

$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(3, 9).Value = 1

$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(3, 10).Value = 1

# Column widths
$ws.Columns.Item(2).ColumnWidth = 17.75
$ws.Columns.Item(4).ColumnWidth = 19.375
$ws.Columns.Item(5).ColumnWidth = 17.75
$ws.Columns.Item(7).ColumnWidth = 8.375
$ws.Columns.Item(8).ColumnWidth = 11.625
$ws.Columns.Item(9).ColumnWidth = 9.625
$ws.Columns.Item(10).ColumnWidth = 10

# Selection moves to A4 (next empty row) as in target
$ws.Range("A4").Select()
